# Menyamakan tabel dengan di laporan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/sheet tab from "barangs" to "barang"
$ws.Name = "barang"

# Move the active selection to C11 (previously C17)
[void]$ws.Range("C11").Select()

# Best-effort: line up the saved window geometry with the reference
# commit (xWindow/yWindow/windowWidth/windowHeight). Harmless no-op if
# the host doesn't persist window chrome geometry.
try {
    $win = $excel.ActiveWindow
    $win.Left = 5400
    $win.Top = 3996
    $win.Width = 13200
    $win.Height = 8964
} catch {
    # window geometry isn't modeled by every host; ignore
}

# Note: D2:D9 hold volatile =RANDBETWEEN(30,60) formulas. Their cached
# <v> results are re-rolled by the engine's recalculation on every
# load/save (by design - RANDBETWEEN is non-deterministic), so the exact
# numbers from the reference commit aren't reproducible nor meaningful to
# hardcode here; the formulas themselves are left untouched.
